# Add a new "GGKP Environment" dataset column (AC) to the Country-Year
# "Datasets and Years" matrix, with checkmarks for every year row (2-29),
# matching the formatting already used by the other dataset columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of column AB (rows 1-29) into the new column AC
$ws.Range("AB1:AB29").Copy()
$ws.Range("AC1:AC29").PasteSpecial(-4122)

# New header for column AC
$ws.Range("AC1").Value = "GGKP Environment"

# Checkmarks for rows 2-29 in column AC
$ws.Range("AC2:AC29").Value = "✓"

# Update the selection / active cell to mirror the saved view state
$null = $ws.Range("AD29").Select()
